# Apply the edits described by the diff to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sistema bloqueado value ---
$ws.Range("C11").Value = "Agua Vapor"

# --- Fechas / horas de inicio y termino ---
$ws.Range("C14").Value = "24/02/2018"
$ws.Range("F14").Value = "2017-04-15T18:16:47Z"
$ws.Range("F15").Value = "2017-04-15T18:16:47Z"
$ws.Range("D20").Value = "2017-04-15T18:16:47Z"
$ws.Range("E20").Value = "2017-04-15T18:16:47Z"
$ws.Range("C15").Value = "25/02/2018"

# --- Itemizado trabajos ---
$ws.Range("C19").Value = "Reparacion estructural"
$ws.Range("C20").Value = "trabajo 1"

# New row: C21 gets the same formatting as C20 (right-aligned label style) plus new text.
$ws.Range("C20").Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4122) | Out-Null
$ws.Range("C21").Value = "trabajo 2"
$excel.CutCopyMode = $false

# --- Observaciones ---
$ws.Range("C35").Value = "Ggggggg"

# Leave the active cell/selection on C21, matching the final saved view state.
$ws.Range("C21").Select() | Out-Null
